# --- 1) Refresh the "time_taken" timestamps on the "data" sheet ---
$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

$data.Range("F2").Value = "2021-10-05 14:33:27.053639"
$data.Range("F3").Value = "2021-10-05 14:33:27.053645"
$data.Range("F4").Value = "2021-10-05 14:33:27.053648"
$data.Range("F5").Value = "2021-10-05 14:33:27.053650"
$data.Range("F6").Value = "2021-10-05 14:33:27.053652"
$data.Range("F7").Value = "2021-10-05 14:33:27.053654"
$data.Range("F8").Value = "2021-10-05 14:33:27.053656"

# --- 2) Add a new "metadata" sheet right after "data" ---
$meta = $wb.Worksheets.Add($null, $data)
$meta.Name = "metadata"

# Header row
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Data row
$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Chondrodysplasia Punctata"
$meta.Range("C2").Value = 70

# data_version ("1.0") must stay text, not be coerced into the number 1
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "1.0"
$meta.Range("D2").Style = "Normal"

$meta.Range("E2").Value = "2021-03-27T06:28:48.005711Z"
$meta.Range("F2").Value = "2021-10-05 14:33:27.051131"
$meta.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/70/?format=json"

# Bold / centered / bordered style, matching the "data" sheet's header style
$headerRange = $meta.Range("B1:G1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$a2 = $meta.Range("A2")
$a2.Font.Bold = $true
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160
$a2.Borders.LineStyle = 1

[void]$meta.Range("A1").Select()

# Keep "data" as the active sheet/selection, as in the original workbook
$data.Activate()
[void]$data.Range("A1").Select()
